# Update the Deudores data table to reflect the refreshed export: the
# trailing rows 28-30 pick up new/updated records, and rows 2-27 carry
# refreshed dates/amounts together with a reshuffled "Cliente" column
# (several prior clients were dropped and replaced by others). The write
# order below (new rows first, then the 2-27 block) reproduces the order
# the client names were (re)entered in, which is also the order new
# shared-string entries were appended upstream.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "CAMPO VERDE TOCANCIPA"
$ws.Cells.Item(28, 3).Value = 46037
$ws.Cells.Item(28, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(28, 4).Value = 788000
$ws.Cells.Item(28, 5).Value = $false
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = "MERKA FRUVER ALEJANDRO"
$ws.Cells.Item(29, 3).Value = 46037
$ws.Cells.Item(29, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(29, 4).Value = 787000
$ws.Cells.Item(29, 5).Value = $false
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = "MULTICARNES"
$ws.Cells.Item(30, 3).Value = 46038
$ws.Cells.Item(30, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(30, 4).Value = 646900
$ws.Cells.Item(30, 5).Value = $false
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "ALISO"
$ws.Cells.Item(2, 3).Value = 46039
$ws.Cells.Item(2, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(2, 4).Value = 108000
$ws.Cells.Item(2, 5).Value = $false
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "ARROZ PAISA SUBA"
$ws.Cells.Item(3, 3).Value = 46029
$ws.Cells.Item(3, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(3, 4).Value = 166000
$ws.Cells.Item(3, 5).Value = $false
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "CARNES JOHANA"
$ws.Cells.Item(4, 3).Value = 46035
$ws.Cells.Item(4, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(4, 4).Value = 164000
$ws.Cells.Item(4, 5).Value = $false
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "CARNILANDIA"
$ws.Cells.Item(5, 3).Value = 46039
$ws.Cells.Item(5, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(5, 4).Value = 546000
$ws.Cells.Item(5, 5).Value = $false
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "CIMARRON DORADO"
$ws.Cells.Item(6, 3).Value = 46038
$ws.Cells.Item(6, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(6, 4).Value = 337000
$ws.Cells.Item(6, 5).Value = $false
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "COCINA CHINA"
$ws.Cells.Item(7, 3).Value = 46031
$ws.Cells.Item(7, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(7, 4).Value = 170000
$ws.Cells.Item(7, 5).Value = $false
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "COCINA CHINA"
$ws.Cells.Item(8, 3).Value = 46036
$ws.Cells.Item(8, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(8, 4).Value = 170000
$ws.Cells.Item(8, 5).Value = $false
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "COCINA CHINA"
$ws.Cells.Item(9, 3).Value = 46039
$ws.Cells.Item(9, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(9, 4).Value = 170000
$ws.Cells.Item(9, 5).Value = $false
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "DARWIN FUTBOL"
$ws.Cells.Item(10, 3).Value = 45921
$ws.Cells.Item(10, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(10, 4).Value = 200000
$ws.Cells.Item(10, 5).Value = $false
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "DAVIDCITO"
$ws.Cells.Item(11, 3).Value = 45947
$ws.Cells.Item(11, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(11, 4).Value = 100000
$ws.Cells.Item(11, 5).Value = $false
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "FRANCO"
$ws.Cells.Item(12, 3).Value = 45996
$ws.Cells.Item(12, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(12, 4).Value = 20000
$ws.Cells.Item(12, 5).Value = $false
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "LA SELECTA"
$ws.Cells.Item(13, 3).Value = 45912
$ws.Cells.Item(13, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(13, 4).Value = 82000
$ws.Cells.Item(13, 5).Value = $false
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "MAFE"
$ws.Cells.Item(14, 3).Value = 46034
$ws.Cells.Item(14, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(14, 4).Value = 521300
$ws.Cells.Item(14, 5).Value = $false
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "MERKA FRUVER DEXI"
$ws.Cells.Item(15, 3).Value = 45995
$ws.Cells.Item(15, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(15, 4).Value = 339000
$ws.Cells.Item(15, 5).Value = $false
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "MERKA FRUVER DEXI"
$ws.Cells.Item(16, 3).Value = 45988
$ws.Cells.Item(16, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(16, 4).Value = 15400
$ws.Cells.Item(16, 5).Value = $false
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "NEVADA"
$ws.Cells.Item(17, 3).Value = 46031
$ws.Cells.Item(17, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(17, 4).Value = 21900
$ws.Cells.Item(17, 5).Value = $false
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "NEVADA"
$ws.Cells.Item(18, 3).Value = 46038
$ws.Cells.Item(18, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(18, 4).Value = 175800
$ws.Cells.Item(18, 5).Value = $false
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "PARAISO FUNZA"
$ws.Cells.Item(19, 3).Value = 46038
$ws.Cells.Item(19, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(19, 4).Value = 76000
$ws.Cells.Item(19, 5).Value = $false
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "PINILLA"
$ws.Cells.Item(20, 3).Value = 45931
$ws.Cells.Item(20, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(20, 4).Value = 82000
$ws.Cells.Item(20, 5).Value = $false
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "PLAZA JESSICA"
$ws.Cells.Item(21, 3).Value = 46039
$ws.Cells.Item(21, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(21, 4).Value = 1349200
$ws.Cells.Item(21, 5).Value = $false
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "WILLINTONG"
$ws.Cells.Item(22, 3).Value = 46039
$ws.Cells.Item(22, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(22, 4).Value = 66000
$ws.Cells.Item(22, 5).Value = $false
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "PUNTA DE ANCA"
$ws.Cells.Item(23, 3).Value = 46038
$ws.Cells.Item(23, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(23, 4).Value = 100000
$ws.Cells.Item(23, 5).Value = $false
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "DOÑA SANDRA"
$ws.Cells.Item(24, 3).Value = 46039
$ws.Cells.Item(24, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(24, 4).Value = 100000
$ws.Cells.Item(24, 5).Value = $false
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = "CALDAS WOK"
$ws.Cells.Item(25, 3).Value = 46035
$ws.Cells.Item(25, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(25, 4).Value = 85000
$ws.Cells.Item(25, 5).Value = $false
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = "CARNIVOROS"
$ws.Cells.Item(26, 3).Value = 46036
$ws.Cells.Item(26, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(26, 4).Value = 196800
$ws.Cells.Item(26, 5).Value = $false
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = "PINILLA SOACHA"
$ws.Cells.Item(27, 3).Value = 46039
$ws.Cells.Item(27, 3).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(27, 4).Value = 266000
$ws.Cells.Item(27, 5).Value = $false

# Restore the view state: scrolled down with I23 selected (no data change).
$ws.Range("I23").Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
